$wb = $excel.ActiveWorkbook

# --- "Data" sheet: the systolic_blood / Emp_status columns (F:G) were random
#     RANDBETWEEN() filler columns that don't belong in the real plots for
#     week-2 exercise, so they get deleted outright (header + data + formulas).
$ws1 = $wb.Worksheets.Item("Data")
$ws1.Columns("F:G").Delete() | Out-Null

# --- "Codebook" sheet: update the two rows that used to document the deleted
#     systolic_blood / Empl_status columns so they instead document the real
#     Chol (Total Cholesterol) and Agecat (Age Chategory) columns.
$ws2 = $wb.Worksheets.Item("Codebook")

$ws2.Range("A5").Value = "Chol"
$ws2.Range("A5").Font.Bold = $false
$ws2.Range("B5").Value = "Total Cholesterol "
$ws2.Range("C5").ClearContents() | Out-Null

$ws2.Range("A6").Value = "Agecat"
$ws2.Range("B6").Value = "Age Chategory "
$ws2.Range("C6").Value = "0 <= 30, 1= 31-40, 2= 41-50, 3=51 or more "

# --- Restore the on-screen selections to match where the author left off ---
$ws2.Range("D6").Select() | Out-Null

$ws1.Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws1.Range("F18").Select() | Out-Null
